$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sVal data (regenerated to filter save games)
$data = @{
    2 = @(3.182878228561681, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 4.173255553662385)
    3 = @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 0, 4.075514443323626)
    4 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054)
    5 = @(0.06328177979961902, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 2.369310727790667)
    6 = @(0.7287194209349384, 1.65323645889881, 16.98373111632243, 6.48142807727062, 0, 25.8471150734268)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 12.0302756157461)
    8 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    9 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G
}

$wb.Save()
